$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(197).Insert()

$ws.Cells.Item(197, 1).Value = 6
$ws.Cells.Item(197, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(197, 3).Value = 'Metropolitana'
$ws.Cells.Item(197, 4).Value = 44798
$ws.Cells.Item(197, 5).Value = 13
$ws.Cells.Item(197, 6).Value = 100112022
$ws.Cells.Item(197, 7).Value = 'Arveja Verde'
$ws.Cells.Item(197, 8).Value = 'Sin especificar'
$ws.Cells.Item(197, 9).Value = 'Primera'
$ws.Cells.Item(197, 10).Value = 250
$ws.Cells.Item(197, 11).Value = 18000
$ws.Cells.Item(197, 12).Value = 20000
$ws.Cells.Item(197, 13).Value = 18800
$ws.Cells.Item(197, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(197, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(197, 16).Value = 752
$ws.Cells.Item(197, 17).Value = 25
$ws.Cells.Item(197, 18).Value = 'Hortaliza'
